$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.708396852016449
$ws.Range("B1").Value = 2.155993938446045
$ws.Range("C1").Value = 3.725496768951416
$ws.Range("D1").Value = 1.174253940582275
$ws.Range("E1").Value = 1.096481204032898
